# Updates the "cryptos" price/volume snapshot in-place to match the
# refreshed GitHub Actions scrape (Mon May 20 18:28:31 UTC 2024).
#
# Notes:
#  - Price (column D) and Volume(1h) (column E) text is refreshed for most
#    rows; two pairs of rows (16/17 and 41/42) also swap their ranking
#    order (Coin name + Link + Price + Volume all move together).
#  - Price values are plain text in the source workbook (e.g. "6.52"),
#    even though they look numeric. Excel auto-converts a bare numeric
#    string typed into a cell into a real number, so for any new price
#    value that parses as a plain number we first force the cell's
#    NumberFormat to Text ("@") to keep it stored as text, matching the
#    original file's representation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.370.52'
$ws.Range("E2").Value = '  +2.51%  '
$ws.Range("D3").Value = '3.143.85'
$ws.Range("E3").Value = '  +2.29%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '578.06'
$ws.Range("E5").Value = '  +0.49%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '180.16'
$ws.Range("E6").Value = '  +6.08%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = '3.142.24'
$ws.Range("E8").Value = '  +2.37%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.522'
$ws.Range("E9").Value = '  +2.12%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.51'
$ws.Range("E10").Value = '  +1.16%  '
$ws.Range("E11").Value = '  +2.00%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.470'
$ws.Range("E12").Value = '  +0.40%  '
$ws.Range("E13").Value = '  +1.68%  '
$ws.Range("E14").Value = '  +3.74%  '
$ws.Range("D15").Value = '3.667.03'
$ws.Range("E15").Value = '  +2.28%  '
$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D16").Value = '68.339.92'
$ws.Range("E16").Value = '  +2.47%  '
$ws.Range("B17").Value = 'TRON'
$ws.Range("C17").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.122'
$ws.Range("E17").Value = '  +0.76%  '
$ws.Range("E18").Value = '  +2.47%  '
$ws.Range("D19").Value = '3.143.93'
$ws.Range("E19").Value = '  +2.59%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.43'
$ws.Range("E20").Value = '  -3.39%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '488.42'
$ws.Range("E21").Value = '  +0.13%  '
$ws.Range("E22").Value = '  +1.87%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.79'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '83.94'
$ws.Range("E25").Value = '  +6.18%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.98'
$ws.Range("E26").Value = '  +2.45%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.58'
$ws.Range("E27").Value = '  +4.03%  '
$ws.Range("E28").Value = '  +0.00%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.09'
$ws.Range("E29").Value = '  +4.15%  '
$ws.Range("E30").Value = '  +4.77%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.64'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '28.26'
$ws.Range("E32").Value = '  +2.50%  '
$ws.Range("E33").Value = '  +0.86%  '
$ws.Range("D34").Value = '0.0₃0950'
$ws.Range("E34").Value = '  +4.45%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").Value = '  +0.00%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '48.42'
$ws.Range("E36").Value = '  +3.02%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.71'
$ws.Range("E37").Value = '  +1.98%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.959'
$ws.Range("E38").Value = '  +1.31%  '
$ws.Range("E39").Value = '  +8.59%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.05'
$ws.Range("E40").Value = '  +4.71%  '
$ws.Range("B41").Value = 'Kaspa'
$ws.Range("C41").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.125'
$ws.Range("E41").Value = '  +2.57%  '
$ws.Range("B42").Value = 'OKB'
$ws.Range("C42").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '49.16'
$ws.Range("E42").Value = '  +0.05%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.41'
$ws.Range("E43").Value = '  +1.73%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.74'
$ws.Range("E44").Value = '  +8.36%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '401.31'
$ws.Range("E45").Value = '  +9.18%  '
$ws.Range("D46").Value = '2.806.49'
$ws.Range("E46").Value = '  +1.84%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '27.54'
$ws.Range("E47").Value = '  +11.47%  '
$ws.Range("E48").Value = '  +1.81%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '135.05'
$ws.Range("E49").Value = '  -0.47%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.38'
$ws.Range("E51").Value = '  +10.85%  '
